$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mislabeled property name that was breaking the file on open:
# "Objekt-ID" -> "ObjektID"
$ws.Range("A2").Value = "ObjektID"

# Build the new (monospace) font used for the data rows on a scratch cell
# first, so the whole font (name + size + colour) is interned as a single
# new style before it gets stamped onto the real range - this mirrors how
# the original fix introduced exactly one new font/cell-style pair.
$scratch = $ws.Range("D1")
$scratch.Font.ThemeColor = 1
$scratch.Font.Name = "Liberation Mono;Courier New;DejaVu Sans Mono;Lucida Sans Typewriter"

$scratch.Copy()
$dataRange = $ws.Range("A2:B32")
$dataRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false
$scratch.Clear()

# Move the active selection back to the top of the sheet (A2) instead of
# leaving it parked on the last cell (B32).
$ws.Range("A2").Select()
